$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update characteristics (Характеристики) column H text for each product row
$ws.Range("H8").Value = "Тип продукта:`nПрофнастил<br>`nОсновной материал:`nСталь<br>`nЦветовая палитра:`nЗеленый<br>`nЦветовая палитра по RAL:`nRAL 6005<br>`nПолезная площадь (м²):`n2.3<br>`nДлина (см):`n200<br>`nШирина (см):`n120<br>`nТолщина (мм):`n0.3<br>`nВысота волны (мм):`n8 <br>`nВес, кг:`n5.99 <br>`nСтрана производства:`nРоссия "
$ws.Range("H7").Value = "Тип продукта:`nПрофнастил <br>`nОсновной материал:`nСталь<br>`nЦветовая палитра:`nКрасный<br>`nЦветовая палитра по RAL:`nRAL 3005<br>`nДлина (см):`n200<br>`nШирина (см):`n120<br>`nТолщина (мм):`n0.3<br>`nВысота волны (мм):`n8<br>`nВес, кг:`n5.99<br>`nСтрана производства:`nРоссия<br>"
$ws.Range("H6").Value = "Тип продукта:`nПрофнастил<br>`nОсновной материал:`nСталь<br>`nЦветовая палитра:`nСерый / серебристый<br>`nЦветовая палитра по RAL:`nНет<br>`nПлощадь покрытия продуктом (м²):`n2.4<br>`nДлина (см):`n200<br>`nШирина (см):`n120<br>`nТолщина (мм):`n0.35<br>`nВысота волны (мм):`n8<br>`nВес на м² (кг):`n3<br>`nВес, кг:`n6.5<br>`nСтрана производства:`nРоссия"
$ws.Range("H4").Value = "Тип продукта:`nСнегозадержатель <br>`nОсновной материал:`nСталь<br>`nПокрытие:`nОцинкованный<br>`nЦветовая палитра по RAL:`nRAL 8017<br>`nЦвет:`nКоричневый<br>`nРазмер (Д х Ш х В) (мм):`n3000х370х170<br>`nВес, кг:`n4.6<br>`nСтрана производства:`nРоссия<br>"
$ws.Range("H3").Value = "Тип продукта:`nПланка карнизная<br>`nОсновной материал:`nПолиэстер<br>`nТолщина (мм):`n0.4<br>`nПокрытие:`nБез покрытия<br>`nЦвет:`nКоричневый<br>`nПрименение продукта:`nЗащищает нижнюю доску обрешётки и лобовую доску от воды<br>`nРазмер (Д х Ш х В) (мм):`n2000 х 200 х 0.4<br>`nВес, кг:`n1.5<br>`nМарка:`nАРТСТРОЙСИТИ<br>`nСтрана производства:`nРоссия<br>"
$ws.Range("H2").Value = "Цвет:`nСерый<br>`nТип продукта:`nМеталлический лист<br>`nВнешний вид поверхности:`nГладкий<br>`nОсновной материал:`nСталь<br>`nДлина (см):`n200<br>`nШирина (см):`n125<br>`nТолщина (мм):`n0.35<br>`nРазмер (Д х Ш х В) (мм):`n2000х1250х0,35<br>`nВес, кг:`n6.5<br>`nПокрытие:`nОцинкованный<br>`nСтрана производства:`nРоссия"

# Fix header style for G1/H1 to match F1 (no-fill variant) by copying format
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore row heights (content edits above trigger autofit; keep original heights)
$ws.Rows.Item(1).RowHeight = 18
$ws.Rows.Item(2).RowHeight = 18
$ws.Rows.Item(3).RowHeight = 17.25
$ws.Rows.Item(4).RowHeight = 17.25
$ws.Rows.Item(5).RowHeight = 17.25
$ws.Rows.Item(6).RowHeight = 17.25
$ws.Rows.Item(7).RowHeight = 17.25
$ws.Rows.Item(8).RowHeight = 17.25

# Update selection to match the recorded cursor position
$ws.Range("G17").Select()

$wb.Save()